$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Globo", "RJ TV 1", "Defesa Civil", "2025-04-07T12:26", "Neutro", "Imagens da praia do Farol de São Tomé. Ondas chegaram a quase 3 metros de altura. Foi feito alerta para não entrar no mar. No sábado, em Farol, 52,2 mm de chuva. Tempo instável. "),
    @("Record", "Balanço Geral", "Saúde", "2025-04-07T12:49", "Positivo", "Começou hoje a campanha de vacinação contra o vírus Influenza em Campos. Repórter *ao vivo*. Muita gente procurando a vacina e consciente da necessidade neste período do ano, quando as temperaturas caem. Expectativa da secretaria de Saúde é vacinar 90% do público alvo. Entrevista com subsecretário de Vigilância em Saúde, Charbell Kury. ")
)

$startRow = 66
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}
